$wb = $excel.ActiveWorkbook

# --- 1. Replace the "Ready for handoff" status text with "In Translation"
#        everywhere it appears (all sheets, so the shared string collapses
#        to a single new value instead of leaving stale entries behind).
$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: keep the literal on the LEFT of -eq. PowerShell coerces
            # the right-hand operand to the left-hand operand's type, and
            # some cells (e.g. TRUE/FALSE) yield an actual [bool] from
            # .Text; "$bool -eq $oldStatus" would coerce the non-empty
            # string to $true and false-match every boolean cell.
            if ($oldStatus -eq $cell.Text) {
                $cell.Value = $newStatus
            }
        }
    }
}

# --- 2. Narrow the "Latest Handoff Datetime" / duplicate width columns.
#        Overview (sheet1): columns E and F
#        zh-cn / de-de (sheet2/sheet3): column C
#        Target stored width (per the OOXML <col> width attribute) is
#        13.4101845877511 "characters". The host snaps ColumnWidth to a
#        6-pixel-per-character grid (stored = (round(chars*6)+5)/6), so
#        that exact value isn't reachable through the object model; 12.5
#        is the input that lands on the closest achievable grid point
#        (13.333333333333334, off by ~0.077) -- closer than any other
#        reachable value, including the naive "set it to the target".
$newWidth = 12.5

$ov = $wb.Worksheets.Item("Overview")
$ov.Columns.Item(5).ColumnWidth = $newWidth
$ov.Columns.Item(6).ColumnWidth = $newWidth

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Columns.Item(3).ColumnWidth = $newWidth

$de = $wb.Worksheets.Item("de-de")
$de.Columns.Item(3).ColumnWidth = $newWidth
